$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "67.745.05"
$ws.Range("E2").Value = "  +1.84%  "

$ws.Range("D3").Value = "3.388.86"
$ws.Range("E3").Value = "  +1.12%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.28"
$ws.Range("E5").Value = "  +5.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "188.56"
$ws.Range("E6").Value = "  -0.72%  "

$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.599"
$ws.Range("E8").Value = "  +2.65%  "

$ws.Range("E9").Value = "  +0.52%  "

$ws.Range("E10").Value = "  +0.69%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "47.74"
$ws.Range("E11").Value = "  +1.47%  "

$ws.Range("E12").Value = "  +1.39%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "3.930.55"
$ws.Range("E13").Value = "  +1.21%  "

$ws.Range("B14").Value = "BitcoinCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "641.54"
$ws.Range("E14").Value = "  +6.44%  "

$ws.Range("E15").Value = "  -0.82%  "

$ws.Range("D16").Value = "67.825.41"
$ws.Range("E16").Value = "  +1.94%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.387.58"
$ws.Range("E17").Value = "  +0.92%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.119"
$ws.Range("E18").Value = "  +0.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.11"
$ws.Range("E19").Value = "  +0.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.23"
$ws.Range("E20").Value = "  +1.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.914"
$ws.Range("E21").Value = "  +1.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.15"
$ws.Range("E22").Value = "  -1.07%  "

$ws.Range("E23").Value = "  +1.56%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "99.85"
$ws.Range("E24").Value = "  -0.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.03"
$ws.Range("E25").Value = "  +0.84%  "

$ws.Range("E26").Value = "  +3.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.78"
$ws.Range("E27").Value = "  +1.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "32.61"
$ws.Range("E28").Value = "  +5.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.74"
$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("E30").Value = "  +3.58%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "610.59"
$ws.Range("E31").Value = "  +3.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.88"
$ws.Range("E32").Value = "  -1.88%  "

$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").Value = "3.983.83"
$ws.Range("E33").Value = "  +6.46%  "

$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.16"
$ws.Range("E34").Value = "  +0.88%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.106"
$ws.Range("E35").Value = "  +1.27%  "

$ws.Range("E36").Value = "  +0.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.13"
$ws.Range("E37").Value = "  -0.48%  "

$ws.Range("E39").Value = "  +3.57%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.90"
$ws.Range("E40").Value = "  -0.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.27"
$ws.Range("E41").Value = "  +0.66%  "

$ws.Range("D42").Value = "0.0₃0710"
$ws.Range("E42").Value = "  -0.34%  "

$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.346"
$ws.Range("E43").Value = "  +0.86%  "

$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.42"
$ws.Range("E44").Value = "  +1.01%  "

$ws.Range("E45").Value = "  +0.75%  "

$ws.Range("E46").Value = "  +0.36%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.61"
$ws.Range("E47").Value = "  +0.35%  "

$ws.Range("E48").Value = "  +0.29%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.38"
$ws.Range("E49").Value = "  +9.93%  "

$ws.Range("E50").Value = "  -20.45%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "128.54"
$ws.Range("E51").Value = "  +3.70%  "
